$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The document window's on-screen position/size at last save (session/display
# artifact of the Excel instance that saved the file). Not meaningfully
# scriptable data, but set here for completeness/best-effort.
$win = $wb.Windows.Item(1)
$win.Left = -98
$win.Top = -98
$win.Width = 21795
$win.Height = 12975

# The user's cursor/selection moved to B3 before the workbook was saved
# (previously I10).
$ws.Range("B3").Select() | Out-Null

# Column C ("Trade Date") is a bestFit/autofit column whose computed width
# grew very slightly (10.42578125 -> 10.7109375 chars) between saves.
# Reproduce via the finest-grained width setter this host exposes.
$ws.Columns.Item(3).ColumnWidth = 9.8
